# Dodanie podziału treningu na części
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Trening" column (F) tagging every sample with the training segment it
# belongs to. Give the header the same look as the other header cells by
# copying A1's format onto it.
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Replace the old 6-row summary with the new, more granular 12-row
# breakdown that also records which part of the training ("Duża Gra" /
# "Mała Gra") each sample belongs to.
$data = @(
    @(45686.47803344907, 1936.4, 10.55, 1.970919694219317, "10-15", "Duża Gra"),
    @(45686.47898599537, 2018.7, 10.11, 1.92757977758135, "10-15", "Duża Gra"),
    @(45686.47945706019, 2059.4, 11.44, 2.198047024863107, "10-15", "Duża Gra"),
    @(45686.47803229166, 1936.3, 9.92, 2.07979691028595, "5-10", "Duża Gra"),
    @(45686.47855081018, 1981.1, 7.8, 2.49568339756557, "5-10", "Duża Gra"),
    @(45686.47898136574, 2018.3, 7.84, 2.738625117710658, "5-10", "Duża Gra"),
    @(45686.48429039352, 2477, 13.6, 3.63463054384504, "10-15", "Mała Gra"),
    @(45686.48802766204, 2799.9, 13.89, 4.436062438147406, "10-15", "Mała Gra"),
    @(45686.48805196759, 2802, 14.3, 2.995907340730942, "10-15", "Mała Gra"),
    @(45686.4842880787, 2476.8, 9.949999999999999, 3.270345296178545, "5-10", "Mała Gra"),
    @(45686.4871931713, 2727.8, 9, 2.836543151310511, "5-10", "Mała Gra"),
    @(45686.48802418981, 2799.6, 8.67, 3.367872629846845, "5-10", "Mała Gra")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}

# Format the timestamp column as a real date-time value. Apply the lowercase
# mask first (registers numFmt 164), then switch to the uppercase mask that
# ends up referenced by the new cell style (numFmt 165); then fan that style
# out to the rest of the timestamp column via copy/paste-format so every row
# shares the one new cell style instead of minting a style per row.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A2").Copy()
$ws.Range("A3:A13").PasteSpecial(-4122)
